$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("#system")

# ------------------------------------------------------------------
# 1) Shift the command columns C:Z one column to the right (-> D:AA)
#    to make room for the new "aws.ses" command column at C.
#    (A plain value copy is used -- not Columns.Insert -- so that the
#    <cols> width metadata, which still refers to raw column indexes,
#    is left completely untouched, matching the source edit.)
# ------------------------------------------------------------------
$src = $ws.Range("C1:Z117")
$dst = $ws.Range("D1")
$src.Copy($dst)

# ------------------------------------------------------------------
# 2) Populate the new column C with the aws.ses command list
# ------------------------------------------------------------------
$ws.Range("C1").Value = "aws.ses"
$ws.Range("C2").Value = "sendMail(profile,to,subject,body)"
$ws.Range("C3").Value = "sendTextMail(profile,to,subject,body)"

# ------------------------------------------------------------------
# 3) Insert "aws.ses" into the target list in column A (sorted
#    between "aws.s3" and "base"), pushing the remaining entries
#    down by one row.
# ------------------------------------------------------------------
$srcA = $ws.Range("A3:A26")
$dstA = $ws.Range("A4")
$srcA.Copy($dstA)
$ws.Range("A3").Value = "aws.ses"

# ------------------------------------------------------------------
# 4) Update the defined names so they continue to point at the
#    correct (now shifted) ranges. A handful of stale/orphaned names
#    (date, db, math, mq, nextgen) are intentionally left untouched,
#    matching the source workbook.
# ------------------------------------------------------------------
$wb.Names.Item("base").RefersTo        = "='#system'!`$D`$2:`$D`$36"
$wb.Names.Item("csv").RefersTo         = "='#system'!`$E`$2:`$E`$5"
$wb.Names.Item("desktop").RefersTo     = "='#system'!`$F`$2:`$F`$92"
$wb.Names.Item("excel").RefersTo       = "='#system'!`$G`$2:`$G`$14"
$wb.Names.Item("external").RefersTo    = "='#system'!`$H`$2:`$H`$3"
$wb.Names.Item("image").RefersTo       = "='#system'!`$I`$2:`$I`$5"
$wb.Names.Item("io").RefersTo          = "='#system'!`$J`$2:`$J`$24"
$wb.Names.Item("jms").RefersTo         = "='#system'!`$K`$2:`$K`$4"
$wb.Names.Item("json").RefersTo        = "='#system'!`$L`$2:`$L`$14"
$wb.Names.Item("mail").RefersTo        = "='#system'!`$M`$2:`$M`$2"
$wb.Names.Item("number").RefersTo      = "='#system'!`$N`$2:`$N`$15"
$wb.Names.Item("pdf").RefersTo         = "='#system'!`$O`$2:`$O`$16"
$wb.Names.Item("rdbms").RefersTo       = "='#system'!`$P`$2:`$P`$7"
$wb.Names.Item("redis").RefersTo       = "='#system'!`$Q`$2:`$Q`$10"
$wb.Names.Item("sms").RefersTo         = "='#system'!`$R`$2:`$R`$2"
$wb.Names.Item("sound").RefersTo       = "='#system'!`$S`$2:`$S`$5"
$wb.Names.Item("ssh").RefersTo         = "='#system'!`$T`$2:`$T`$9"
$wb.Names.Item("step").RefersTo        = "='#system'!`$U`$2:`$U`$4"
$wb.Names.Item("target").RefersTo      = "='#system'!`$A`$2:`$A`$27"
$wb.Names.Item("web").RefersTo         = "='#system'!`$V`$2:`$V`$117"
$wb.Names.Item("webalert").RefersTo    = "='#system'!`$W`$2:`$W`$8"
$wb.Names.Item("webcookie").RefersTo   = "='#system'!`$X`$2:`$X`$8"
$wb.Names.Item("ws").RefersTo          = "='#system'!`$Y`$2:`$Y`$17"
$wb.Names.Item("ws.async").RefersTo    = "='#system'!`$Z`$2:`$Z`$8"
$wb.Names.Item("xml").RefersTo         = "='#system'!`$AA`$2:`$AA`$11"

$wb.Names.Add("aws.ses", "='#system'!`$C`$2:`$C`$3")
